# Apply dataset cleanup edit:
#  - remove the placeholder hyperlinks on column K (site URLs were all
#    pointing at http://www.google.com with rich-text hyperlink formatting)
#  - replace placeholder "Lorem" occupation text with each candidate's
#    real occupation
#  - replace placeholder site URLs with each candidate's real site URL
#    (Kenton Johnson has none, so it stays blank)
#  - replace the generic "Photo courtesy of the candidate" photo-credit
#    text with the real photo credit per candidate (blank where no photo
#    credit applies)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all hyperlinks from the sheet (column K used to hyperlink to
# http://www.google.com as filler/placeholder links).
$ws.Hyperlinks.Delete()

# --- Row 2: Aspen Dunaway ---
$ws.Range("C2").Value = "Attorney"
$ws.Range("K2").Value = "http://www.dunawayfortexas.com"
$ws.Range("L2").Value = "Rachel Zein"

# --- Row 3: Huey Rey Fischer ---
$ws.Range("C3").Value = "Former legislative aide"
$ws.Range("K3").Value = "http://hueyfischer.com"
$ws.Range("L3").Value = "Rachel Zein"

# --- Row 4: Gina Hinojosa ---
$ws.Range("C4").Value = "Austin ISD School Board member"
$ws.Range("K4").Value = "http://www.ginaforaustin.com"
$ws.Range("L4").Value = "Rachel Zein"

# --- Row 5: Kenton Johnson ---
$ws.Range("C5").Value = "Attorney"
$ws.Range("K5").ClearContents()
$ws.Range("L5").ClearContents()

# --- Row 6: Blake Rocap ---
$ws.Range("C6").Value = "Former NARAL Pro-Choice Texas counsel"
$ws.Range("K6").Value = "http://www.blakerocap.com"
$ws.Range("L6").Value = "Blake Rocap campaign"

# --- Row 7: Matt Shrum ---
$ws.Range("C7").Value = "Attorney"
$ws.Range("K7").Value = "http://www.mattshrumfortexas.com"
$ws.Range("L7").ClearContents()

# --- Row 8: Heather Way ---
$ws.Range("C8").Value = "University of Texas law professor"
$ws.Range("K8").Value = "http://www.voteheatherway.com"
$ws.Range("L8").Value = "Heather Way campaign"
